$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 144, shifting existing rows 144:161 down to 145:162
$ws.Rows.Item(144).Insert()

# Populate the newly inserted row 144 with the new record's data
$ws.Cells.Item(144, 1).Value = 7
$ws.Cells.Item(144, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(144, 3).Value = "Ñuble"
$ws.Cells.Item(144, 4).Value = 44491
$ws.Cells.Item(144, 5).Value = 16
$ws.Cells.Item(144, 6).Value = 100112043
$ws.Cells.Item(144, 7).Value = "Pepino ensalada"
$ws.Cells.Item(144, 8).Value = "Sin especificar"
$ws.Cells.Item(144, 9).Value = "Primera"
$ws.Cells.Item(144, 10).Value = 120
$ws.Cells.Item(144, 11).Value = 13000
$ws.Cells.Item(144, 12).Value = 14000
$ws.Cells.Item(144, 13).Value = 13500
$ws.Cells.Item(144, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(144, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(144, 16).Value = 225
$ws.Cells.Item(144, 17).Value = 60
$ws.Cells.Item(144, 18).Value = "Hortaliza"
